$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3734236.8
$ws.Range("I132").Value = 3734236.8
$ws.Range("K132").Value = 11202710.4
$ws.Range("M132").Value = -11200180.4
$ws.Range("H137").Value = 986327.5
$ws.Range("I137").Value = 1365.1111
$ws.Range("J137").Value = 2168282.2
$ws.Range("K137").Value = 4095.3333
$ws.Range("L137").Value = 6504846.600000001
$ws.Range("M137").Value = -1545.3333
$ws.Range("N137").Value = -6509946.600000001
$ws.Range("H138").Value = 4336.5117
$ws.Range("I138").Value = 1698.1333
$ws.Range("J138").Value = 5749.9287
$ws.Range("K138").Value = 5094.3999
$ws.Range("L138").Value = 17249.7861
$ws.Range("M138").Value = 45.60009999999966
$ws.Range("N138").Value = -27529.7861
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2641
$ws.Range("I32").Value = 2710.5938
$ws.Range("K32").Value = 2710.5938
$ws.Range("M32").Value = -2423.5938
$ws.Range("H61").Value = 1669751.9
$ws.Range("I61").Value = 2224409.2
$ws.Range("K61").Value = 2224409.2
$ws.Range("M61").Value = -2224197.2
$ws.Range("H74").Value = 3138.5862
$ws.Range("I74").Value = 1484.8948
$ws.Range("J74").Value = 6280.6
$ws.Range("K74").Value = 1484.8948
$ws.Range("L74").Value = 6280.6
$ws.Range("M74").Value = -610.8948
$ws.Range("N74").Value = -8028.6
$ws.Range("H77").Value = 3138.5862
$ws.Range("I77").Value = 1484.8948
$ws.Range("J77").Value = 6280.6
$ws.Range("K77").Value = 7424.474
$ws.Range("L77").Value = 31403
$ws.Range("M77").Value = -3056.474
$ws.Range("N77").Value = -40139
$ws.Range("H97").Value = 1647
$ws.Range("I97").Value = 1623.8
$ws.Range("J97").Value = 1763
$ws.Range("K97").Value = 1623.8
$ws.Range("L97").Value = 1763
$ws.Range("M97").Value = -1127.8
$ws.Range("N97").Value = -2755
$ws.Range("H113").Value = 100318.2
$ws.Range("J113").Value = 100318.2
$ws.Range("L113").Value = 100318.2
$ws.Range("N113").Value = -108996.2
$ws.Range("H132").Value = 410161.5
$ws.Range("I132").Value = 466765.66
$ws.Range("K132").Value = 1400296.98
$ws.Range("M132").Value = -1397766.98
$ws.Range("H136").Value = 1669751.9
$ws.Range("I136").Value = 2224409.2
$ws.Range("K136").Value = 6673227.600000001
$ws.Range("M136").Value = -6670677.600000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2065.389
$ws.Range("I86").Value = 1993.2858
$ws.Range("J86").Value = 2317.75
$ws.Range("K86").Value = 1993.2858
$ws.Range("L86").Value = 2317.75
$ws.Range("M86").Value = -870.2858000000001
$ws.Range("N86").Value = -4563.75
$ws.Range("H89").Value = 2065.389
$ws.Range("I89").Value = 1993.2858
$ws.Range("J89").Value = 2317.75
$ws.Range("K89").Value = 9966.429
$ws.Range("L89").Value = 11588.75
$ws.Range("M89").Value = -4350.429
$ws.Range("N89").Value = -22820.75
$ws.Range("H134").Value = 2485383.2
$ws.Range("I134").Value = 3062370.5
$ws.Range("K134").Value = 9187111.5
$ws.Range("M134").Value = -9184576.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 237731.98
$ws.Range("I31").Value = 355484.88
$ws.Range("J31").Value = 20818.736
$ws.Range("K31").Value = 355484.88
$ws.Range("L31").Value = 20818.736
$ws.Range("M31").Value = -355189.88
$ws.Range("N31").Value = -21408.736
$ws.Range("H34").Value = 237731.98
$ws.Range("I34").Value = 355484.88
$ws.Range("J34").Value = 20818.736
$ws.Range("K34").Value = 355484.88
$ws.Range("L34").Value = 20818.736
$ws.Range("M34").Value = -355282.88
$ws.Range("N34").Value = -21222.736
$ws.Range("H41").Value = 38009.832
$ws.Range("J41").Value = 50000
$ws.Range("L41").Value = 50000
$ws.Range("N41").Value = -50856
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H134").Value = 5620.7036
$ws.Range("I134").Value = 6400.0435
$ws.Range("J134").Value = 1139.5
$ws.Range("K134").Value = 19200.1305
$ws.Range("L134").Value = 3418.5
$ws.Range("M134").Value = -16665.1305
$ws.Range("N134").Value = -8488.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 128.6
$ws.Range("I2").Value = 42
$ws.Range("K2").Value = 252
$ws.Range("M2").Value = -139
$ws.Range("H17").Value = 1028.1111
$ws.Range("I17").Value = 1316.1428
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 3948.4284
$ws.Range("L17").Value = 60
$ws.Range("M17").Value = -3779.4284
$ws.Range("N17").Value = -398
$ws.Range("H55").Value = 2494
$ws.Range("I55").Value = 737
$ws.Range("J55").Value = 4836.6665
$ws.Range("K55").Value = 2211
$ws.Range("L55").Value = 14509.9995
$ws.Range("M55").Value = -2034
$ws.Range("N55").Value = -14863.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2474.6858
$ws.Range("I102").Value = 1496.7391
$ws.Range("J102").Value = 4349.0835
$ws.Range("K102").Value = 1496.7391
$ws.Range("L102").Value = 4349.0835
$ws.Range("M102").Value = 125.2609
$ws.Range("N102").Value = -7593.0835
$ws.Range("H126").Value = 796948.1
$ws.Range("I126").Value = 1668906.2
$ws.Range("J126").Value = 4259
$ws.Range("K126").Value = 5006718.6
$ws.Range("L126").Value = 12777
$ws.Range("M126").Value = -5004248.6
$ws.Range("N126").Value = -17717
$ws.Range("H132").Value = 356931.38
$ws.Range("I132").Value = 432934.56
$ws.Range("J132").Value = 2249.8333
$ws.Range("K132").Value = 1298803.68
$ws.Range("L132").Value = 6749.499899999999
$ws.Range("M132").Value = -1296273.68
$ws.Range("N132").Value = -11809.4999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8982.263000000001
$ws.Range("J16").Value = 3817.3333
$ws.Range("L16").Value = 3817.3333
$ws.Range("N16").Value = -4157.3333
$ws.Range("H40").Value = 3825.3125
$ws.Range("I40").Value = 3469.3845
$ws.Range("K40").Value = 3469.3845
$ws.Range("M40").Value = -3333.3845
$ws.Range("H42").Value = 20025
$ws.Range("I42").Value = 20025
$ws.Range("K42").Value = 20025
$ws.Range("M42").Value = -19462
$ws.Range("H46").Value = 2176.4119
$ws.Range("J46").Value = 2607.5
$ws.Range("L46").Value = 2607.5
$ws.Range("N46").Value = -2983.5
$ws.Range("H49").Value = 20025
$ws.Range("I49").Value = 20025
$ws.Range("K49").Value = 20025
$ws.Range("M49").Value = -19878
$ws.Range("H93").Value = 1890.7
$ws.Range("I93").Value = 1738.375
$ws.Range("K93").Value = 1738.375
$ws.Range("M93").Value = -490.375
$ws.Range("H122").Value = 4690.909
$ws.Range("I122").Value = 4514.5713
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 13543.7139
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -11093.7139
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 4354725
$ws.Range("I132").Value = 4354725
$ws.Range("K132").Value = 13064175
$ws.Range("M132").Value = -13061645
$ws.Range("H136").Value = 1958.8
$ws.Range("I136").Value = 2136
$ws.Range("J136").Value = 1250
$ws.Range("K136").Value = 6408
$ws.Range("L136").Value = 3750
$ws.Range("M136").Value = -3858
$ws.Range("N136").Value = -8850
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 66666.664
$ws.Range("J46").Value = 66666.664
$ws.Range("L46").Value = 66666.664
$ws.Range("N46").Value = -67128.664
$ws.Range("H122").Value = 1566.4445
$ws.Range("I122").Value = 1566.4445
$ws.Range("K122").Value = 4699.333500000001
$ws.Range("M122").Value = -2249.333500000001
$ws.Range("H132").Value = 7190768
$ws.Range("I132").Value = 10594484
$ws.Range("K132").Value = 31783452
$ws.Range("M132").Value = -31780922
$ws.Range("H134").Value = 66666.664
$ws.Range("J134").Value = 66666.664
$ws.Range("L134").Value = 199999.992
$ws.Range("N134").Value = -205069.992
$ws.Range("H136").Value = 71556.5
$ws.Range("I136").Value = 71556.5
$ws.Range("K136").Value = 214669.5
$ws.Range("M136").Value = -212119.5
